# Bombátor now loads CSV files instead of Excel files, so the previous
# workbook's duplicated "EDAV" rows (the result of the old per-teacher-id
# Excel-based loading) are collapsed back down to a single row.
#
# Rows 7-14 (1-based worksheet rows) are exact duplicates of row 6
# ("EDAV"/"Data Analysis and Visualisation") and are removed entirely,
# shifting rows 15-26 up to become rows 7-18. Excel automatically shrinks
# the table ("Frame0") range, its AutoFilter, and the sheet dimension from
# A1:F26 down to A1:F18 as part of the row delete.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7:A14").EntireRow.Delete()
